$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr2d = New-Object 'object[,]' 16,14
$arr2d[0,0] = 54.53585066666667
$arr2d[0,1] = 163.607552
$arr2d[0,2] = 0.3031388658437607
$arr2d[0,3] = 0.3031388658437607
$arr2d[0,4] = 3
$arr2d[0,5] = 1
$arr2d[0,6] = 0.2577576666666667
$arr2d[0,7] = 0.7732730000000001
$arr2d[0,8] = 0.03524815007985697
$arr2d[0,9] = 0.03524815007985697
$arr2d[0,10] = 14.057033617521782
$arr2d[0,11] = 126.51330255769601
$arr2d[0,12] = 0.010685084238298505
$arr2d[0,13] = 0.010685084238298505
$arr2d[1,0] = 54.53585066666667
$arr2d[1,1] = 163.607552
$arr2d[1,2] = 0.3031388658437607
$arr2d[1,3] = 0.3031388658437607
$arr2d[1,4] = 3
$arr2d[1,5] = 1
$arr2d[1,6] = 2.242708666666667
$arr2d[1,7] = 6.728126
$arr2d[1,8] = 0.3066885757089511
$arr2d[1,9] = 0.3066885757089511
$arr2d[1,10] = 122.30802493417248
$arr2d[1,11] = 1100.772224407552
$arr2d[1,12] = 0.09296922700764979
$arr2d[1,13] = 0.09296922700764976
$arr2d[2,0] = 54.53585066666667
$arr2d[2,1] = 163.607552
$arr2d[2,2] = 0.3031388658437607
$arr2d[2,3] = 0.3031388658437607
$arr2d[2,4] = 3
$arr2d[2,5] = 1
$arr2d[2,6] = 4.789377333333333
$arr2d[2,7] = 14.368132
$arr2d[2,8] = 0.6549434327891901
$arr2d[2,9] = 0.6549434327891901
$arr2d[2,10] = 261.1927670369849
$arr2d[2,11] = 2350.7349033328637
$arr2d[2,12] = 0.1985388094075344
$arr2d[2,13] = 0.1985388094075344
$arr2d[3,0] = 54.53585066666667
$arr2d[3,1] = 163.607552
$arr2d[3,2] = 0.3031388658437607
$arr2d[3,3] = 0.3031388658437607
$arr2d[3,4] = 1
$arr2d[3,5] = 0.3333333333333333
$arr2d[3,6] = 0.02281433333333334
$arr2d[3,7] = 0.068443
$arr2d[3,8] = 0.0031198414220018683
$arr2d[3,9] = 0.0031198414220018683
$arr2d[3,10] = 1.2441990757262227
$arr2d[3,11] = 11.197791681536
$arr2d[3,12] = 0.000945745190278032
$arr2d[3,13] = 0.0009457451902780318
$arr2d[4,0] = 18.46467533333333
$arr2d[4,1] = 55.394026
$arr2d[4,2] = 0.10263635150631549
$arr2d[4,3] = 0.10263635150631549
$arr2d[4,4] = 3
$arr2d[4,5] = 1
$arr2d[4,6] = 0.2577576666666667
$arr2d[4,7] = 0.7732730000000001
$arr2d[4,8] = 0.03524815007985697
$arr2d[4,9] = 0.03524815007985697
$arr2d[4,10] = 4.7594116296775555
$arr2d[4,11] = 42.834704667098
$arr2d[4,12] = 0.003617741521543562
$arr2d[4,13] = 0.003617741521543563
$arr2d[5,0] = 18.46467533333333
$arr2d[5,1] = 55.394026
$arr2d[5,2] = 0.10263635150631549
$arr2d[5,3] = 0.10263635150631549
$arr2d[5,4] = 3
$arr2d[5,5] = 1
$arr2d[5,6] = 2.242708666666667
$arr2d[5,7] = 6.728126
$arr2d[5,8] = 0.3066885757089511
$arr2d[5,9] = 0.3066885757089511
$arr2d[5,10] = 41.410887397252885
$arr2d[5,11] = 372.69798657527593
$arr2d[5,12] = 0.031477396459435154
$arr2d[5,13] = 0.031477396459435154
$arr2d[6,0] = 18.46467533333333
$arr2d[6,1] = 55.394026
$arr2d[6,2] = 0.10263635150631549
$arr2d[6,3] = 0.10263635150631549
$arr2d[6,4] = 3
$arr2d[6,5] = 1
$arr2d[6,6] = 4.789377333333333
$arr2d[6,7] = 14.368132
$arr2d[6,8] = 0.6549434327891901
$arr2d[6,9] = 0.6549434327891901
$arr2d[6,10] = 88.43429750882575
$arr2d[6,11] = 795.9086775794319
$arr2d[6,12] = 0.06722100438450422
$arr2d[6,13] = 0.06722100438450423
$arr2d[7,0] = 18.46467533333333
$arr2d[7,1] = 55.394026
$arr2d[7,2] = 0.10263635150631549
$arr2d[7,3] = 0.10263635150631549
$arr2d[7,4] = 1
$arr2d[7,5] = 0.3333333333333333
$arr2d[7,6] = 0.02281433333333334
$arr2d[7,7] = 0.068443
$arr2d[7,8] = 0.0031198414220018683
$arr2d[7,9] = 0.0031198414220018683
$arr2d[7,10] = 0.4212592579464444
$arr2d[7,11] = 3.791333321518
$arr2d[7,12] = 0.00032020914083254686
$arr2d[7,13] = 0.0003202091408325469
$arr2d[8,0] = 12.55635966666667
$arr2d[8,1] = 37.669079
$arr2d[8,2] = 0.06979483370938175
$arr2d[8,3] = 0.06979483370938175
$arr2d[8,4] = 3
$arr2d[8,5] = 1
$arr2d[8,6] = 0.2577576666666667
$arr2d[8,7] = 0.7732730000000001
$arr2d[8,8] = 0.03524815007985697
$arr2d[8,9] = 0.03524815007985697
$arr2d[8,10] = 3.2364979695074463
$arr2d[8,11] = 29.128481725567006
$arr2d[8,12] = 0.0024601387733869484
$arr2d[8,13] = 0.002460138773386948
$arr2d[9,0] = 12.55635966666667
$arr2d[9,1] = 37.669079
$arr2d[9,2] = 0.06979483370938175
$arr2d[9,3] = 0.06979483370938175
$arr2d[9,4] = 3
$arr2d[9,5] = 1
$arr2d[9,6] = 2.242708666666667
$arr2d[9,7] = 6.728126
$arr2d[9,8] = 0.3066885757089511
$arr2d[9,9] = 0.3066885757089511
$arr2d[9,10] = 28.160256646217125
$arr2d[9,11] = 253.44230981595402
$arr2d[9,12] = 0.021405278142173377
$arr2d[9,13] = 0.02140527814217337
$arr2d[10,0] = 12.55635966666667
$arr2d[10,1] = 37.669079
$arr2d[10,2] = 0.06979483370938175
$arr2d[10,3] = 0.06979483370938175
$arr2d[10,4] = 3
$arr2d[10,5] = 1
$arr2d[10,6] = 4.789377333333333
$arr2d[10,7] = 14.368132
$arr2d[10,8] = 0.6549434327891901
$arr2d[10,9] = 0.6549434327891901
$arr2d[10,10] = 60.13714437671424
$arr2d[10,11] = 541.2342993904281
$arr2d[10,12] = 0.04571166798057316
$arr2d[10,13] = 0.04571166798057316
$arr2d[11,0] = 12.55635966666667
$arr2d[11,1] = 37.669079
$arr2d[11,2] = 0.06979483370938175
$arr2d[11,3] = 0.06979483370938175
$arr2d[11,4] = 1
$arr2d[11,5] = 0.3333333333333333
$arr2d[11,6] = 0.02281433333333334
$arr2d[11,7] = 0.068443
$arr2d[11,8] = 0.0031198414220018683
$arr2d[11,9] = 0.0031198414220018683
$arr2d[11,10] = 0.2864649748885557
$arr2d[11,11] = 2.5781847739970005
$arr2d[11,12] = 0.0002177488132482615
$arr2d[11,13] = 0.00021774881324826143
$arr2d[12,0] = 94.34696966666667
$arr2d[12,1] = 283.040909
$arr2d[12,2] = 0.5244299489405421
$arr2d[12,3] = 0.5244299489405421
$arr2d[12,4] = 3
$arr2d[12,5] = 1
$arr2d[12,6] = 0.2577576666666667
$arr2d[12,7] = 0.7732730000000001
$arr2d[12,8] = 0.03524815007985697
$arr2d[12,9] = 0.03524815007985697
$arr2d[12,10] = 24.31865475835078
$arr2d[12,11] = 218.867892825157
$arr2d[12,12] = 0.018485185546627952
$arr2d[12,13] = 0.018485185546627956
$arr2d[13,0] = 94.34696966666667
$arr2d[13,1] = 283.040909
$arr2d[13,2] = 0.5244299489405421
$arr2d[13,3] = 0.5244299489405421
$arr2d[13,4] = 3
$arr2d[13,5] = 1
$arr2d[13,6] = 2.242708666666667
$arr2d[13,7] = 6.728126
$arr2d[13,8] = 0.3066885757089511
$arr2d[13,9] = 0.3066885757089511
$arr2d[13,10] = 211.5927665451705
$arr2d[13,11] = 1904.334898906534
$arr2d[13,12] = 0.1608366740996928
$arr2d[13,13] = 0.1608366740996928
$arr2d[14,0] = 94.34696966666667
$arr2d[14,1] = 283.040909
$arr2d[14,2] = 0.5244299489405421
$arr2d[14,3] = 0.5244299489405421
$arr2d[14,4] = 3
$arr2d[14,5] = 1
$arr2d[14,6] = 4.789377333333333
$arr2d[14,7] = 14.368132
$arr2d[14,8] = 0.6549434327891901
$arr2d[14,9] = 0.6549434327891901
$arr2d[14,10] = 451.8632379902209
$arr2d[14,11] = 4066.769141911988
$arr2d[14,12] = 0.3434719510165783
$arr2d[14,13] = 0.34347195101657835
$arr2d[15,0] = 94.34696966666667
$arr2d[15,1] = 283.040909
$arr2d[15,2] = 0.5244299489405421
$arr2d[15,3] = 0.5244299489405421
$arr2d[15,4] = 1
$arr2d[15,5] = 0.3333333333333333
$arr2d[15,6] = 0.02281433333333334
$arr2d[15,7] = 0.068443
$arr2d[15,8] = 0.0031198414220018683
$arr2d[15,9] = 0.0031198414220018683
$arr2d[15,10] = 2.152463214965223
$arr2d[15,11] = 19.372168934687
$arr2d[15,12] = 0.001636138277643028
$arr2d[15,13] = 0.0016361382776430278

$ws.Range("G2:T17").Value = $arr2d
